$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 98
$ws.Range("H98").Value = 326567.03
$ws.Range("I98").Value = 4068.4783
$ws.Range("J98").Value = 1253750.4
$ws.Range("K98").Value = 4068.4783
$ws.Range("L98").Value = 1253750.4
$ws.Range("M98").Value = -2570.4783
$ws.Range("N98").Value = -1256746.4
# Row 113
$ws.Range("H113").Value = 2252.4707
$ws.Range("I113").Value = 2018
$ws.Range("J113").Value = 2587.4285
$ws.Range("K113").Value = 2018
$ws.Range("L113").Value = 2587.4285
$ws.Range("M113").Value = 1236
$ws.Range("N113").Value = -9095.4285
# Row 122
$ws.Range("H122").Value = 326567.03
$ws.Range("I122").Value = 4068.4783
$ws.Range("J122").Value = 1253750.4
$ws.Range("K122").Value = 12205.4349
$ws.Range("L122").Value = 3761251.2
$ws.Range("M122").Value = -9755.4349
$ws.Range("N122").Value = -3766151.2
# Row 129
$ws.Range("H129").Value = 880.14703
$ws.Range("I129").Value = 397.57144
$ws.Range("J129").Value = 1005.2593
$ws.Range("K129").Value = 1192.71432
$ws.Range("L129").Value = 3015.7779
$ws.Range("M129").Value = 3807.28568
$ws.Range("N129").Value = -13015.7779

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 400224
$ws.Range("I32").Value = 2786.5522
$ws.Range("K32").Value = 2786.5522
$ws.Range("M32").Value = -2499.5522

$ws = $wb.Worksheets.Item("BSM")
# Row 80
$ws.Range("H80").Value = 502.10526
$ws.Range("J80").Value = 392.16666
$ws.Range("L80").Value = 392.16666
$ws.Range("N80").Value = -2388.16666
# Row 83
$ws.Range("H83").Value = 502.10526
$ws.Range("J83").Value = 392.16666
$ws.Range("L83").Value = 1960.8333
$ws.Range("N83").Value = -11944.8333
# Row 107
$ws.Range("H107").Value = 732.26666
$ws.Range("I107").Value = 541.6786
$ws.Range("J107").Value = 3400.5
$ws.Range("K107").Value = 541.6786
$ws.Range("L107").Value = 3400.5
$ws.Range("M107").Value = 1378.3214
$ws.Range("N107").Value = -7240.5
# Row 140
$ws.Range("H140").Value = 59461.54
$ws.Range("J140").Value = 59461.54
$ws.Range("L140").Value = 59461.54
$ws.Range("N140").Value = -69821.54000000001

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 7464599
$ws.Range("I31").Value = 9805494
$ws.Range("J31").Value = 2996.0625
$ws.Range("K31").Value = 9805494
$ws.Range("L31").Value = 2996.0625
$ws.Range("M31").Value = -9805199
$ws.Range("N31").Value = -3586.0625
# Row 34
$ws.Range("H34").Value = 7464599
$ws.Range("I34").Value = 9805494
$ws.Range("J34").Value = 2996.0625
$ws.Range("K34").Value = 9805494
$ws.Range("L34").Value = 2996.0625
$ws.Range("M34").Value = -9805292
$ws.Range("N34").Value = -3400.0625
# Row 99
$ws.Range("H99").Value = 1491.7273
$ws.Range("I99").Value = 1490.9
$ws.Range("K99").Value = 1490.9
$ws.Range("M99").Value = 7.099999999999909
# Row 122
$ws.Range("H122").Value = 1578.4445
$ws.Range("I122").Value = 1555.6364
$ws.Range("J122").Value = 1614.2858
$ws.Range("K122").Value = 4666.9092
$ws.Range("L122").Value = 4842.857400000001
$ws.Range("M122").Value = -2216.9092
$ws.Range("N122").Value = -9742.857400000001
# Row 126
$ws.Range("H126").Value = 1491.7273
$ws.Range("I126").Value = 1490.9
$ws.Range("K126").Value = 4472.700000000001
$ws.Range("M126").Value = -2002.700000000001

$ws = $wb.Worksheets.Item("CUL")
# Row 14
$ws.Range("H14").Value = 122.818184
$ws.Range("I14").Value = 122.818184
$ws.Range("K14").Value = 368.454552
$ws.Range("M14").Value = -195.454552
# Row 131
$ws.Range("H131").Value = 9260053
$ws.Range("J131").Value = 11111974
$ws.Range("L131").Value = 33335922
$ws.Range("N131").Value = -33346002

$ws = $wb.Worksheets.Item("GSM")
# Row 14
$ws.Range("H14").Value = 20000002
$ws.Range("I14").Value = 20000002
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 20000002
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -19999834
$ws.Range("N14").ClearContents()
# Row 70
$ws.Range("H70").Value = 15635017
$ws.Range("I70").Value = 40192820
$ws.Range("J70").Value = 7323.364
$ws.Range("K70").Value = 40192820
$ws.Range("L70").Value = 7323.364
$ws.Range("M70").Value = -40192550
$ws.Range("N70").Value = -7863.364
# Row 73
$ws.Range("H73").Value = 15635017
$ws.Range("I73").Value = 40192820
$ws.Range("J73").Value = 7323.364
$ws.Range("K73").Value = 40192820
$ws.Range("L73").Value = 7323.364
$ws.Range("M73").Value = -40191884
$ws.Range("N73").Value = -9195.364
# Row 122
$ws.Range("H122").Value = 2277.1304
$ws.Range("I122").Value = 1972.3158
$ws.Range("K122").Value = 5916.9474
$ws.Range("M122").Value = -3466.9474
# Row 126
$ws.Range("H126").Value = 8334799
$ws.Range("I126").Value = 1333.1111
$ws.Range("J126").Value = 15153089
$ws.Range("K126").Value = 3999.3333
$ws.Range("L126").Value = 45459267
$ws.Range("M126").Value = -1529.3333
$ws.Range("N126").Value = -45464207

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3460.889
$ws.Range("I7").Value = 2921.7778
$ws.Range("J7").Value = 4000
$ws.Range("K7").Value = 2921.7778
$ws.Range("L7").Value = 4000
$ws.Range("M7").Value = -2809.7778
$ws.Range("N7").Value = -4224
# Row 126
$ws.Range("H126").Value = 3460.889
$ws.Range("I126").Value = 2921.7778
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 8765.3334
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -6295.3334
$ws.Range("N126").Value = -16940

$ws = $wb.Worksheets.Item("WVR")
# Row 37
$ws.Range("H37").Value = 3299.4
$ws.Range("I37").Value = 2850
$ws.Range("K37").Value = 2850
$ws.Range("M37").Value = -2647
# Row 41
$ws.Range("H41").Value = 6576.1665
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 6576.1665
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 6576.1665
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -7356.1665
# Row 74
$ws.Range("H74").Value = 12820.833
$ws.Range("J74").Value = 12820.833
$ws.Range("L74").Value = 12820.833
$ws.Range("N74").Value = -14692.833
# Row 77
$ws.Range("H77").Value = 12820.833
$ws.Range("J77").Value = 12820.833
$ws.Range("L77").Value = 38462.499
$ws.Range("N77").Value = -47822.499
# Row 126
$ws.Range("H126").Value = 1680.5807
$ws.Range("I126").Value = 1230.2
$ws.Range("J126").Value = 2499.4546
$ws.Range("K126").Value = 3690.6
$ws.Range("L126").Value = 7498.3638
$ws.Range("M126").Value = -1220.6
$ws.Range("N126").Value = -12438.3638

Write-Output "done"